# Add "Sweet Home 3D" entry to the Chocolatey package list on Tabelle1.
# This inserts a new row at position 100 (shifting the existing rows 100-125
# down to 101-126, and 127-129 down to 128-130), sets the new row's data
# ("sweet-home-3d" in column A, "Grundriss, Raumplaner" comment in column G),
# and keeps the sheet's autofilter / filter-database defined name in sync
# with the now-larger data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert a new blank row above the current row 100 (sysinternals), shifting
# everything below it down by one row.
$ws.Rows(100).Insert(-4121)   # -4121 = xlShiftDown

# The freshly inserted row picks up a blended/merged style when created via
# Insert(); copy the (now shifted-down) formatting from the row that used to
# be row 100 so the new row matches the existing table's look exactly.
$ws.Range("A101:G101").Copy()
$ws.Range("A100:G100").PasteSpecial(-4122)  # -4122 = xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new entry's data.
$ws.Range("A100").Value = "sweet-home-3d"
$ws.Range("G100").Value = "Grundriss, Raumplaner"

# Re-apply the autofilter so its range grows to cover the new last row.
$ws.AutoFilterMode = $false
$ws.Range("A1:G126").AutoFilter(1)

# Keep the hidden _FilterDatabase defined name (used by the autofilter) in
# sync with the new range as well.
$fdName = $wb.Names.Item("Tabelle1!_FilterDatabase")
$fdName.RefersTo = "=Tabelle1!`$A`$1:`$G`$126"
